$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Add the new "pastTimeHorizon" parameter row below the existing table
$ws.Range("A12").Value = "pastTimeHorizon"
$ws.Range("A12").WrapText = $true
$ws.Range("B12").Value = 5

# Move the active selection as recorded in the saved workbook
$ws.Range("D18").Select()
